$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.666571333333333
$ws.Range("H2").Value = 4.999714
$ws.Range("I2").Value = 0.6861303711299348
$ws.Range("J2").Value = 0.6861303711299347
$ws.Range("Q2").Value = 0.01058328349044444
$ws.Range("R2").Value = 0.09524955141399999
$ws.Range("S2").Value = 0.6861303711299348
$ws.Range("T2").Value = 0.6861303711299347

# Row 3
$ws.Range("I3").Value = 0.2417614358401214
$ws.Range("J3").Value = 0.2417614358401214
$ws.Range("S3").Value = 0.2417614358401214
$ws.Range("T3").Value = 0.2417614358401214

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1751466666666667
$ws.Range("H4").Value = 0.52544
$ws.Range("I4").Value = 0.0721081930299439
$ws.Range("J4").Value = 0.0721081930299439
$ws.Range("Q4").Value = 0.001112239715555555
$ws.Range("R4").Value = 0.01001015744
$ws.Range("S4").Value = 0.0721081930299439
$ws.Range("T4").Value = 0.0721081930299439
